# Re-style the three tables (slides 14, 15, 16) that used the deck's
# default table style so they use the built-in "Medium Style 2 - Accent 1"
# table style instead.
#
#   {7F0AF831-DB1A-4BA0-8FEF-5FFD96BE88F1}  ->  {26449846-B909-4E95-9BFD-53AF5B09280A}

$p = $ppt.ActivePresentation

$slideIndexes = 14, 15, 16
$newStyleId = "{26449846-B909-4E95-9BFD-53AF5B09280A}"

foreach ($idx in $slideIndexes) {
    $slide = $p.Slides.Item($idx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId)
        }
    }
}
